# Updates the cryptos price-list sheet to the latest scraped snapshot.
# Column D ("Price") cells store numeric-looking values as TEXT in the
# source data (e.g. "599.25", "0.0000272", "1.00"), so each write is
# forced to text with a leading apostrophe to avoid Excel's automatic
# "looks like a number" coercion (which would otherwise drop things like
# trailing zeros or switch to scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.972.63"
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = "'3.212.91"
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'602.48"
$ws.Range('E5').Value = '  +4.53%  '
$ws.Range('D6').Value = "'151.91"
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'3.211.82"
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').Value = "'0.538"
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('D11').Value = "'6.14"
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = "'0.512"
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = "'0.0000272"
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = "'38.66"
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = "'3.731.26"
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = "'66.012.84"
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = "'7.41"
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').Value = "'3.210.80"
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = "'0.112"
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = "'512.93"
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = "'15.81"
$ws.Range('E21').Value = '  +5.55%  '
$ws.Range('D22').Value = "'0.739"
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = "'15.23"
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').Value = "'7.99"
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value = "'85.44"
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'9.29"
$ws.Range('E27').Value = '  +2.11%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'3.03"
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').Value = "'2.25"
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').Value = "'2.89"
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = "'6.82"
$ws.Range('E31').Value = '  +7.53%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'28.20"
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = "'6.66"
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').Value = "'55.52"
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').Value = "'0.0926"
$ws.Range('E37').Value = '  +3.34%  '
$ws.Range('D38').Value = "'488.77"
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('D39').Value = "'0.0425"
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = "'3.01"
$ws.Range('E40').Value = '  -3.54%  '
$ws.Range('D41').Value = "'8.91"
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('D42').Value = "'3.028.22"
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').Value = "'0.120"
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').Value = "'0.293"
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('E45').Value = '  +5.85%  '
$ws.Range('D46').Value = "'2.47"
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').Value = "'29.21"
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = "'0.117"
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = "'2.33"
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').Value = "'120.53"
$ws.Range('E51').Value = '  -1.35%  '
